$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.6950908042689878"
$ws.Range("E2").Value = [double]"0.6950908042689878"

$ws.Range("D3").Value = [double]"0.001418715245402583"
$ws.Range("E3").Value = [double]"0.001418715245402583"

$ws.Range("D4").Value = [double]"0.9710666448946109"
$ws.Range("E4").Value = [double]"0.9710666448946109"

$ws.Range("D5").Value = [double]"0.001136821777982218"
$ws.Range("E5").Value = [double]"0.001136821777982218"

$ws.Range("D6").Value = [double]"0.1911117487871187"
$ws.Range("E6").Value = [double]"0.1911117487871187"

$ws.Range("D7").Value = [double]"0.9999999999683244"
$ws.Range("E7").Value = [double]"3.167555107097542E-11"

$ws.Range("D8").Value = [double]"0.9582423534428004"
$ws.Range("E8").Value = [double]"0.0417576465571996"

$ws.Range("D9").Value = [double]"0.761090458840944"
$ws.Range("E9").Value = [double]"0.238909541159056"

$ws.Range("D10").Value = [double]"0.5053742700161901"
$ws.Range("E10").Value = [double]"0.4946257299838099"

$ws.Range("C11").Value = $false
$ws.Range("D11").Value = [double]"1.178463594580837E-05"
$ws.Range("E11").Value = [double]"0.9999882153640542"
$ws.Range("F11").Value = [double]"1.729198098182678"
$ws.Range("G11").Value = [double]"0.7"

$ws.Range("D12").Value = [double]"0.8865139484026694"
$ws.Range("E12").Value = [double]"0.8865139484026694"

$ws.Range("D13").Value = [double]"0.01346969077979133"
$ws.Range("E13").Value = [double]"0.01346969077979133"

$ws.Range("D14").Value = [double]"0.9860184529618947"
$ws.Range("E14").Value = [double]"0.9860184529618947"

$ws.Range("D15").Value = [double]"0.008569857483065163"
$ws.Range("E15").Value = [double]"0.008569857483065163"

$ws.Range("D16").Value = [double]"0.1251051171349855"
$ws.Range("E16").Value = [double]"0.1251051171349855"

$ws.Range("D17").Value = [double]"0.9999999999999998"
$ws.Range("E17").Value = [double]"2.220446049250313E-16"

$ws.Range("D18").Value = [double]"0.9974594307080419"
$ws.Range("E18").Value = [double]"0.002540569291958117"

$ws.Range("D19").Value = [double]"0.7637299225550545"
$ws.Range("E19").Value = [double]"0.2362700774449455"

$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"0.05033173115500862"
$ws.Range("E20").Value = [double]"0.9496682688449913"

$ws.Range("C21").Value = $false
$ws.Range("D21").Value = [double]"0.1253193904637937"
$ws.Range("E21").Value = [double]"0.8746806095362063"
$ws.Range("F21").Value = [double]"1.194000482559204"
$ws.Range("G21").Value = [double]"0.6"
